# Add a new column E to the worksheet with header value 3 and
# per-row fractional values, matching the style of column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1 - copy style from D1 (bold/bordered header style)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E1").Value = 3

# Data values for E2:E12
$values = @(
    0.001565,
    0.001129,
    0.000321,
    0.000224,
    0.00019,
    0.000157,
    0.000158,
    0.00015,
    0.000535,
    0.002207,
    0.000041
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
